$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns (F and G) before the existing lat/long columns.
#    This shifts the previous F (위도/lat) and G (경도/long) columns to H/I.
# ---------------------------------------------------------------------------
$ws.Columns("F:G").Insert()

# ---------------------------------------------------------------------------
# 2. Fix up the formatting (style) of the newly inserted cells.
#    Columns("F:G").Insert() copies the format of the column to the left
#    (column E) for each row, which is not what we want for every row, so
#    explicitly re-apply the correct formats:
#      - header row (row 1)  -> same style as the rest of the header (A1 etc.)
#      - row 2                -> same style as the rest of row 2 (A2 etc.)
#      - rows 3-37            -> same style as the rest of those rows (A3 etc.)
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("F2:G2").PasteSpecial(-4122) | Out-Null

$ws.Range("A3").Copy() | Out-Null
$ws.Range("F3:G37").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Set the column widths for the two new columns.
# ---------------------------------------------------------------------------
$ws.Range("F:F").ColumnWidth = 8.571428571428571
$ws.Range("G:G").ColumnWidth = 4.857142857142857

# ---------------------------------------------------------------------------
# 4. Fill in the header labels and the values for the new columns.
#    F = 준공년도 (completion year), G = 지역 (region)
#    (values are entered in this particular order so that the resulting
#    shared-string table lines up with the authored workbook)
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "지역"
$ws.Range("F1").Value = "준공년도"

$ws.Range("G2:G37").Value = "수도권"
$ws.Range("F2:F37").Value = "2024년"

# ---------------------------------------------------------------------------
# 5. Update the sheet view: zoom to 85%, no fixed top-left cell, and update
#    the active selection.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Range("F3:F37").Select() | Out-Null
